# "Fruta / hortaliza, semanal" — weekly data refresh.
#
# A new weekly price record (fecha serial 44476 = 2021-10-07) is inserted
# at row 160, pushing the existing rows 160..193 down to 161..194 (the
# sheet grows from A1:R193 to A1:R194). Everything below the new row keeps
# its original values - only its row position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 160, shifting row 160 (and everything below
# it) down by one. Excel's default insert behaviour (shift down) mirrors
# the formatting of the row above, which is what carries the date-style
# (s="2") onto the new D160 cell automatically.
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(160, 1).Value = 5
$ws.Cells.Item(160, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(160, 3).Value = "Maule"
$ws.Cells.Item(160, 4).Value = 44476
$ws.Cells.Item(160, 5).Value = 7
$ws.Cells.Item(160, 6).Value = 100114013
$ws.Cells.Item(160, 7).Value = "Zanahoria"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 400
$ws.Cells.Item(160, 11).Value = 8000
$ws.Cells.Item(160, 12).Value = 8000
$ws.Cells.Item(160, 13).Value = 8000
$ws.Cells.Item(160, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(160, 15).Value = "Región de Ñuble"
$ws.Cells.Item(160, 16).Value = 400
$ws.Cells.Item(160, 17).Value = 20
$ws.Cells.Item(160, 18).Value = "Hortaliza"
